$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores numeric-looking values as text (e.g. "0.529",
# "2.411.32"). Force a Text number format on the whole data range so Excel
# keeps the updated cells as literal strings instead of coercing them to
# numbers (which would drop trailing zeros / collapse the dotted look).
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "62.114.77"
$ws.Range("E2").Value = "  -0.22%  "

# Row 3
$ws.Range("D3").Value = "2.411.32"
$ws.Range("E3").Value = "  -0.59%  "

# Row 4
$ws.Range("E4").Value = "  -0.18%  "

# Row 5
$ws.Range("D5").Value = "562.41"
$ws.Range("E5").Value = "  +1.22%  "

# Row 6
$ws.Range("D6").Value = "142.46"
$ws.Range("E6").Value = "  -0.81%  "

# Row 7
$ws.Range("E7").Value = "  +0.08%  "

# Row 8
$ws.Range("D8").Value = "0.529"
$ws.Range("E8").Value = "  -0.44%  "

# Row 9
$ws.Range("D9").Value = "0.109"
$ws.Range("E9").Value = "  +0.02%  "

# Row 10
$ws.Range("E10").Value = "  -2.08%  "

# Row 11
$ws.Range("D11").Value = "5.29"
$ws.Range("E11").Value = "  -1.92%  "

# Row 12
$ws.Range("D12").Value = "0.350"
$ws.Range("E12").Value = "  -1.57%  "

# Row 13
$ws.Range("D13").Value = "25.55"
$ws.Range("E13").Value = "  -2.92%  "

# Row 14
$ws.Range("D14").Value = "0.0000173"
$ws.Range("E14").Value = "  -1.16%  "

# Row 15
$ws.Range("D15").Value = "2.846.74"
$ws.Range("E15").Value = "  -0.65%  "

# Row 16
$ws.Range("D16").Value = "62.090.76"
$ws.Range("E16").Value = "  +0.58%  "

# Row 17
$ws.Range("D17").Value = "2.406.65"
$ws.Range("E17").Value = "  -1.04%  "

# Row 18
$ws.Range("D18").Value = "11.28"
$ws.Range("E18").Value = "  +0.79%  "

# Row 19
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "4.14"
$ws.Range("E19").Value = "  -1.40%  "

# Row 20
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "320.75"
$ws.Range("E20").Value = "  -1.36%  "

# Row 21
$ws.Range("D21").Value = "6.83"
$ws.Range("E21").Value = "  +0.53%  "

# Row 22
$ws.Range("E22").Value = "  -0.03%  "

# Row 23
$ws.Range("D23").Value = "66.06"
$ws.Range("E23").Value = "  +1.41%  "

# Row 24
$ws.Range("D24").Value = "1.72"
$ws.Range("E24").Value = "  -1.80%  "

# Row 25
$ws.Range("D25").Value = "8.85"
$ws.Range("E25").Value = "  -2.84%  "

# Row 26
$ws.Range("D26").Value = "569.77"
$ws.Range("E26").Value = "  +1.26%  "

# Row 27
$ws.Range("E27").Value = "  +0.67%  "

# Row 28
$ws.Range("D28").Value = "2.528.07"
$ws.Range("E28").Value = "  +1.01%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0940"
$ws.Range("E29").Value = "  +0.22%  "

# Row 30
$ws.Range("D30").Value = "8.17"
$ws.Range("E30").Value = "  -2.04%  "

# Row 31
$ws.Range("D31").Value = "1.41"
$ws.Range("E31").Value = "  -3.12%  "

# Row 32
$ws.Range("D32").Value = "0.148"
$ws.Range("E32").Value = "  -0.81%  "

# Row 33
$ws.Range("E33").Value = "  -0.17%  "

# Row 34
$ws.Range("D34").Value = "1.52"
$ws.Range("E34").Value = "  -2.52%  "

# Row 35
$ws.Range("E35").Value = "  +0.12%  "

# Row 36
$ws.Range("D36").Value = "4.68"
$ws.Range("E36").Value = "  -2.94%  "

# Row 37
$ws.Range("D37").Value = "5.47"
$ws.Range("E37").Value = "  -6.02%  "

# Row 38
$ws.Range("D38").Value = "0.380"
$ws.Range("E38").Value = "  -1.31%  "

# Row 39
$ws.Range("D39").Value = "151.41"
$ws.Range("E39").Value = "  +3.31%  "

# Row 40
$ws.Range("D40").Value = "18.61"
$ws.Range("E40").Value = "  -1.15%  "

# Row 41
$ws.Range("D41").Value = "1.79"
$ws.Range("E41").Value = "  -9.53%  "

# Row 43
$ws.Range("D43").Value = "2.28"
$ws.Range("E43").Value = "  -0.06%  "

# Row 44
$ws.Range("D44").Value = "147.68"
$ws.Range("E44").Value = "  -1.26%  "

# Row 45
$ws.Range("D45").Value = "3.62"
$ws.Range("E45").Value = "  -0.72%  "

# Row 46
$ws.Range("D46").Value = "0.0531"
$ws.Range("E46").Value = "  -1.99%  "

# Row 47
$ws.Range("D47").Value = "19.84"
$ws.Range("E47").Value = "  -2.69%  "

# Row 48
$ws.Range("D48").Value = "0.592"
$ws.Range("E48").Value = "  -0.32%  "

# Row 49
$ws.Range("D49").Value = "0.0917"
$ws.Range("E49").Value = "  +0.31%  "

# Row 50
$ws.Range("D50").Value = "0.0225"
$ws.Range("E50").Value = "  -0.94%  "

# Row 51
$ws.Range("D51").Value = "1.06"
$ws.Range("E51").Value = "  +4.39%  "
